$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1: add a new paragraph "Git rebase -abort = Cancela el rebase
# actual." right after the "Git rebase [branch] = Manda los cambios..."
# paragraph, moving the _GoBack bookmark from the end of the old
# paragraph to the end of the new one (matches the diff's hunk #1).
# -----------------------------------------------------------------------

$enDash = [char]0x2013

$findRng = $d.Content
$findRng.Find.Execute("especificado en el comando.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$findRng.Collapse(0)  # wdCollapseEnd -> lands right before the existing _GoBack bookmark

# Insert the paragraph break first; this pushes the (degenerate) _GoBack
# bookmark into the newly created (still empty) paragraph.
$findRng.InsertBefore([char]13)

# Insert the new paragraph's text immediately before the bookmark's
# current position (in order) -- each insertion keeps landing just to the
# left of the bookmark, so when done in forward order the bookmark ends
# up sitting right after all of the new text, i.e. at the end of the new
# paragraph, matching the target document.
$bm = $d.Bookmarks.Item("_GoBack")
$insRng = $d.Range($bm.Start, $bm.Start)
$insRng.InsertBefore("Git rebase " + $enDash)

$bm = $d.Bookmarks.Item("_GoBack")
$insRng = $d.Range($bm.Start, $bm.Start)
$insRng.InsertBefore("abort")

$bm = $d.Bookmarks.Item("_GoBack")
$insRng = $d.Range($bm.Start, $bm.Start)
$insRng.InsertBefore(" = Cancela el rebase actual.")

# -----------------------------------------------------------------------
# Change 2: add an extra empty paragraph right before the
# "Git commit -m [mensaje] = ..." paragraph (matches the diff's hunk #2).
# -----------------------------------------------------------------------

$commitIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Git commit -m *") {
        $commitIndex = $i
        break
    }
}

$prevPara = $d.Paragraphs.Item($commitIndex - 1)
$prevPara.Range.InsertParagraphAfter()

Write-Output "Done"
